# error solve ifrs list
# Corrects the financial figures (rows 2-9, columns D:AJ) on the
# "company_list" sheet for 팬오션 (Pan Ocean) - the previous values were
# off by roughly two orders of magnitude / wrong units for most metrics.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 16456
$ws.Cells.Item(2, 5).Value = 2160
$ws.Cells.Item(2, 6).Value = 2160
$ws.Cells.Item(2, 7).Value = 7901
$ws.Cells.Item(2, 8).Value = 7861
$ws.Cells.Item(2, 9).Value = 7861
$ws.Cells.Item(2, 10).Value = -1
$ws.Cells.Item(2, 11).Value = 45579
$ws.Cells.Item(2, 12).Value = 31353
$ws.Cells.Item(2, 13).Value = 14226
$ws.Cells.Item(2, 14).Value = 14213
$ws.Cells.Item(2, 15).Value = 13
$ws.Cells.Item(2, 16).Value = 2145
$ws.Cells.Item(2, 17).Value = 3874
$ws.Cells.Item(2, 18).Value = 3065
$ws.Cells.Item(2, 19).Value = -5034
$ws.Cells.Item(2, 20).Value = 417
$ws.Cells.Item(2, 21).Value = 3457
$ws.Cells.Item(2, 22).Value = 17668
$ws.Cells.Item(2, 23).Value = 13.12
$ws.Cells.Item(2, 24).Value = 47.77
$ws.Cells.Item(2, 25).Value = 95.06
$ws.Cells.Item(2, 26).Value = 16.66
$ws.Cells.Item(2, 27).Value = 220.4
$ws.Cells.Item(2, 28).Value = 454.81
$ws.Cells.Item(2, 29).Value = 6060
$ws.Cells.Item(2, 30).Value = 0.54
$ws.Cells.Item(2, 31).Value = 8301
$ws.Cells.Item(2, 32).Value = 0.39
$ws.Cells.Item(2, 33).Value = 0
$ws.Cells.Item(2, 34).Value = 0
$ws.Cells.Item(2, 35).Value = 0
$ws.Cells.Item(2, 36).Value = 171210697

# Row 3
$ws.Cells.Item(3, 4).Value = 18193
$ws.Cells.Item(3, 5).Value = 2294
$ws.Cells.Item(3, 6).Value = 2294
$ws.Cells.Item(3, 7).Value = 478
$ws.Cells.Item(3, 8).Value = 455
$ws.Cells.Item(3, 9).Value = 455
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 43143
$ws.Cells.Item(3, 12).Value = 18828
$ws.Cells.Item(3, 13).Value = 24314
$ws.Cells.Item(3, 14).Value = 24062
$ws.Cells.Item(3, 15).Value = 253
$ws.Cells.Item(3, 16).Value = 5244
$ws.Cells.Item(3, 17).Value = 3053
$ws.Cells.Item(3, 18).Value = -392
$ws.Cells.Item(3, 19).Value = -4171
$ws.Cells.Item(3, 20).Value = 547
$ws.Cells.Item(3, 21).Value = 2506
$ws.Cells.Item(3, 22).Value = 15199
$ws.Cells.Item(3, 23).Value = 12.61
$ws.Cells.Item(3, 24).Value = 2.5
$ws.Cells.Item(3, 25).Value = 2.38
$ws.Cells.Item(3, 26).Value = 1.03
$ws.Cells.Item(3, 27).Value = 77.44
$ws.Cells.Item(3, 28).Value = 306.08
$ws.Cells.Item(3, 29).Value = 125
$ws.Cells.Item(3, 30).Value = 29.62
$ws.Cells.Item(3, 31).Value = 4533
$ws.Cells.Item(3, 32).Value = 0.8100000000000001
$ws.Cells.Item(3, 33).Value = 0
$ws.Cells.Item(3, 34).Value = 0
$ws.Cells.Item(3, 35).Value = 0
$ws.Cells.Item(3, 36).Value = 530758755

# Row 4
$ws.Cells.Item(4, 4).Value = 18740
$ws.Cells.Item(4, 5).Value = 1679
$ws.Cells.Item(4, 6).Value = 1679
$ws.Cells.Item(4, 7).Value = 991
$ws.Cells.Item(4, 8).Value = 971
$ws.Cells.Item(4, 9).Value = 979
$ws.Cells.Item(4, 10).Value = -8
$ws.Cells.Item(4, 11).Value = 43306
$ws.Cells.Item(4, 12).Value = 17648
$ws.Cells.Item(4, 13).Value = 25659
$ws.Cells.Item(4, 14).Value = 25407
$ws.Cells.Item(4, 15).Value = 252
$ws.Cells.Item(4, 16).Value = 5344
$ws.Cells.Item(4, 17).Value = 2510
$ws.Cells.Item(4, 18).Value = -924
$ws.Cells.Item(4, 19).Value = -1979
$ws.Cells.Item(4, 20).Value = 994
$ws.Cells.Item(4, 21).Value = 1516
$ws.Cells.Item(4, 22).Value = 15203
$ws.Cells.Item(4, 23).Value = 8.960000000000001
$ws.Cells.Item(4, 24).Value = 5.18
$ws.Cells.Item(4, 25).Value = 3.96
$ws.Cells.Item(4, 26).Value = 2.25
$ws.Cells.Item(4, 27).Value = 68.78
$ws.Cells.Item(4, 28).Value = 323.87
$ws.Cells.Item(4, 29).Value = 184
$ws.Cells.Item(4, 30).Value = 21.66
$ws.Cells.Item(4, 31).Value = 4754
$ws.Cells.Item(4, 32).Value = 0.84
$ws.Cells.Item(4, 33).Value = 0
$ws.Cells.Item(4, 34).Value = 0
$ws.Cells.Item(4, 35).Value = 0
$ws.Cells.Item(4, 36).Value = 534396771

# Row 5
$ws.Cells.Item(5, 4).Value = 23362
$ws.Cells.Item(5, 5).Value = 1950
$ws.Cells.Item(5, 6).Value = 1950
$ws.Cells.Item(5, 7).Value = 1432
$ws.Cells.Item(5, 8).Value = 1413
$ws.Cells.Item(5, 9).Value = 1431
$ws.Cells.Item(5, 10).Value = -18
$ws.Cells.Item(5, 11).Value = 38944
$ws.Cells.Item(5, 12).Value = 14846
$ws.Cells.Item(5, 13).Value = 24098
$ws.Cells.Item(5, 14).Value = 23890
$ws.Cells.Item(5, 15).Value = 207
$ws.Cells.Item(5, 16).Value = 5345
$ws.Cells.Item(5, 17).Value = 2616
$ws.Cells.Item(5, 18).Value = -1172
$ws.Cells.Item(5, 19).Value = -1550
$ws.Cells.Item(5, 20).Value = 1065
$ws.Cells.Item(5, 21).Value = 1551
$ws.Cells.Item(5, 22).Value = 12221
$ws.Cells.Item(5, 23).Value = 8.35
$ws.Cells.Item(5, 24).Value = 6.05
$ws.Cells.Item(5, 25).Value = 5.8
$ws.Cells.Item(5, 26).Value = 3.44
$ws.Cells.Item(5, 27).Value = 61.61
$ws.Cells.Item(5, 28).Value = 350.78
$ws.Cells.Item(5, 29).Value = 268
$ws.Cells.Item(5, 30).Value = 19.69
$ws.Cells.Item(5, 31).Value = 4469
$ws.Cells.Item(5, 32).Value = 1.18
$ws.Cells.Item(5, 33).Value = 0
$ws.Cells.Item(5, 34).Value = 0
$ws.Cells.Item(5, 35).Value = 0
$ws.Cells.Item(5, 36).Value = 534537812

# Row 6
$ws.Cells.Item(6, 4).Value = 26684
$ws.Cells.Item(6, 5).Value = 2039
$ws.Cells.Item(6, 6).Value = 2039
$ws.Cells.Item(6, 7).Value = 1490
$ws.Cells.Item(6, 8).Value = 1486
$ws.Cells.Item(6, 9).Value = 1524
$ws.Cells.Item(6, 11).Value = 41195
$ws.Cells.Item(6, 12).Value = 14591
$ws.Cells.Item(6, 13).Value = 26604
$ws.Cells.Item(6, 14).Value = 26438
$ws.Cells.Item(6, 16).Value = 5346
$ws.Cells.Item(6, 17).Value = 2681
$ws.Cells.Item(6, 18).Value = -1266
$ws.Cells.Item(6, 19).Value = -1491
$ws.Cells.Item(6, 20).Value = 1654
$ws.Cells.Item(6, 21).Value = 1026
$ws.Cells.Item(6, 22).Value = 11270
$ws.Cells.Item(6, 23).Value = 7.64
$ws.Cells.Item(6, 24).Value = 5.57
$ws.Cells.Item(6, 25).Value = 6.06
$ws.Cells.Item(6, 26).Value = 3.71
$ws.Cells.Item(6, 27).Value = 54.84
$ws.Cells.Item(6, 28).Value = 379.25
$ws.Cells.Item(6, 29).Value = 285
$ws.Cells.Item(6, 30).Value = 15.59
$ws.Cells.Item(6, 31).Value = 4946
$ws.Cells.Item(6, 32).Value = 0.9
$ws.Cells.Item(6, 33).ClearContents()
$ws.Cells.Item(6, 34).ClearContents()
$ws.Cells.Item(6, 35).Value = 0
$ws.Cells.Item(6, 36).Value = 534569207

# Row 7
$ws.Cells.Item(7, 4).Value = 24947
$ws.Cells.Item(7, 5).Value = 2165
$ws.Cells.Item(7, 7).Value = 1631
$ws.Cells.Item(7, 8).Value = 1617
$ws.Cells.Item(7, 9).Value = 1647
$ws.Cells.Item(7, 11).Value = 44179
$ws.Cells.Item(7, 12).Value = 14734
$ws.Cells.Item(7, 13).Value = 29445
$ws.Cells.Item(7, 14).Value = 29246
$ws.Cells.Item(7, 16).Value = 5346
$ws.Cells.Item(7, 17).Value = 3499
$ws.Cells.Item(7, 18).Value = -2546
$ws.Cells.Item(7, 19).Value = -778
$ws.Cells.Item(7, 20).Value = 2779
$ws.Cells.Item(7, 21).Value = 1325
$ws.Cells.Item(7, 23).Value = 8.68
$ws.Cells.Item(7, 24).Value = 6.48
$ws.Cells.Item(7, 25).Value = 5.92
$ws.Cells.Item(7, 26).Value = 3.79
$ws.Cells.Item(7, 27).Value = 50.04
$ws.Cells.Item(7, 29).Value = 308
$ws.Cells.Item(7, 30).Value = 12.92
$ws.Cells.Item(7, 31).Value = 5471
$ws.Cells.Item(7, 32).Value = 0.73
$ws.Cells.Item(7, 33).Value = 6
$ws.Cells.Item(7, 34).Value = 0.15
$ws.Cells.Item(7, 35).Value = 1.89

# Row 8
$ws.Cells.Item(8, 4).Value = 28827
$ws.Cells.Item(8, 5).Value = 2474
$ws.Cells.Item(8, 7).Value = 1927
$ws.Cells.Item(8, 8).Value = 1873
$ws.Cells.Item(8, 9).Value = 1892
$ws.Cells.Item(8, 11).Value = 46543
$ws.Cells.Item(8, 12).Value = 14949
$ws.Cells.Item(8, 13).Value = 31595
$ws.Cells.Item(8, 14).Value = 31427
$ws.Cells.Item(8, 16).Value = 5346
$ws.Cells.Item(8, 17).Value = 3547
$ws.Cells.Item(8, 18).Value = -2380
$ws.Cells.Item(8, 19).Value = -444
$ws.Cells.Item(8, 20).Value = 2509
$ws.Cells.Item(8, 21).Value = 1200
$ws.Cells.Item(8, 23).Value = 8.58
$ws.Cells.Item(8, 24).Value = 6.5
$ws.Cells.Item(8, 25).Value = 6.24
$ws.Cells.Item(8, 26).Value = 4.13
$ws.Cells.Item(8, 27).Value = 47.31
$ws.Cells.Item(8, 29).Value = 354
$ws.Cells.Item(8, 30).Value = 11.25
$ws.Cells.Item(8, 31).Value = 5879
$ws.Cells.Item(8, 32).Value = 0.68
$ws.Cells.Item(8, 33).Value = 10
$ws.Cells.Item(8, 34).Value = 0.25
$ws.Cells.Item(8, 35).Value = 2.83

# Row 9
$ws.Cells.Item(9, 4).Value = 30533
$ws.Cells.Item(9, 5).Value = 2696
$ws.Cells.Item(9, 7).Value = 2155
$ws.Cells.Item(9, 8).Value = 2067
$ws.Cells.Item(9, 9).Value = 2085
$ws.Cells.Item(9, 11).Value = 48338
$ws.Cells.Item(9, 12).Value = 14411
$ws.Cells.Item(9, 13).Value = 33928
$ws.Cells.Item(9, 14).Value = 33792
$ws.Cells.Item(9, 16).Value = 5346
$ws.Cells.Item(9, 17).Value = 3852
$ws.Cells.Item(9, 18).Value = -1938
$ws.Cells.Item(9, 19).Value = -573
$ws.Cells.Item(9, 20).Value = 2112
$ws.Cells.Item(9, 21).Value = 1746
$ws.Cells.Item(9, 23).Value = 8.83
$ws.Cells.Item(9, 24).Value = 6.77
$ws.Cells.Item(9, 25).Value = 6.39
$ws.Cells.Item(9, 26).Value = 4.36
$ws.Cells.Item(9, 27).Value = 42.47
$ws.Cells.Item(9, 29).Value = 390
$ws.Cells.Item(9, 30).Value = 10.2
$ws.Cells.Item(9, 31).Value = 6321
$ws.Cells.Item(9, 32).Value = 0.63
$ws.Cells.Item(9, 33).Value = 12
$ws.Cells.Item(9, 34).Value = 0.29
$ws.Cells.Item(9, 35).Value = 2.99
